# Populate the three new bank-detail columns (Q/R/S) added to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Cells.Item(1, 17).Value = "bank"
$ws.Cells.Item(1, 18).Value = "bankAccName"
$ws.Cells.Item(1, 19).Value = "bankAccount"

# Per-row bank number / bank account name / bank account number
$bankNum = @(1, 2, 1, 2, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31)
$bankAccName = @("altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana", "nasanjargal", "altanbagana")
$bankAccount = @(1232155, 54674567, 108116979, 161559391, 215001803, 268444215, 321886627, 375329039, 428771451, 482213863, 535656275, 589098687, 642541099, 695983511, 749425923, 802868335, 856310747, 909753159, 963195571, 1016637983, 1070080395, 1123522807, 1176965219, 1230407631, 1283850043, 1337292455, 1390734867, 1444177279, 1497619691, 1551062103, 1604504515)

for ($i = 0; $i -lt $bankNum.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $bankNum[$i]
    $ws.Cells.Item($row, 18).Value = $bankAccName[$i]
    $ws.Cells.Item($row, 19).Value = $bankAccount[$i]
}

# Widen column G and size the two new columns (closest widths this host can represent)
$ws.Columns.Item(7).ColumnWidth = 5.833333333333334
$ws.Columns.Item(18).ColumnWidth = 12.666666666666668
$ws.Columns.Item(19).ColumnWidth = 11.666666666666668

# Move the active selection to T5, matching the saved view state
$ws.Range("T5").Select() | Out-Null
